# Generate Report for Handoff
# Updates the "Latest Handoff/Handback" timestamp for the 99c3d242 row (row 5 / row index
# 4 of the data, i.e. spreadsheet row 5) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-10-13 12:56:41"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2016-10-13 12:56:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2016-10-13 12:56:41"
